$d = $word.ActiveDocument

# --- 1. Update the ID placeholder text in the first paragraph --------------
# The paragraph currently reads:
#   "**ID__AFFARS_5301_topic_9__ID**" + " "   (two runs, identical rPr)
# It must become just:
#   "**ID__AFFARS_SUBPART_5301_3__ID**"       (a single run, no trailing space)

$oldId = "**ID__AFFARS_5301_topic_9__ID**"
$newId = "**ID__AFFARS_SUBPART_5301_3__ID**"

# Locate the placeholder text without touching anything yet.
$idRange = $d.Content.Duplicate
$found = $idRange.Find.Execute($oldId, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Replace just that run's text in place (keeps it as its own run, distinct
# from the following space run).
$idRange.Text = $newId

# The space run now immediately follows the new text; locate and remove it.
$spaceRange = $d.Range($idRange.End, $idRange.End + 1)
if ($spaceRange.Text -eq " ") {
    $spaceRange.Text = ""
}

# --- 2. Paragraph formatting on that same (first) paragraph ----------------
$p1 = $d.Paragraphs(1)

# Add a paragraph border (top/left/bottom/right) with 5pt space-from-text
# and no explicit line (matches <w:pBdr><w:top w:space="5"/>...).
$p1.Borders.DistanceFromTop = 5
$p1.Borders.DistanceFromLeft = 5
$p1.Borders.DistanceFromBottom = 5
$p1.Borders.DistanceFromRight = 5

# Left indent: 120 twips (6pt) -> 225 twips (11.25pt)
$p1.Range.ParagraphFormat.LeftIndent = 11.25

Write-Output "done"
